$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'61.781.02"
$ws.Range("E2").Value = "  +1.43%  "
$ws.Range("D3").Value = "'3.452.72"
$ws.Range("E3").Value = "  +3.61%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.26%  "
$ws.Range("D5").Value = "'579.71"
$ws.Range("E5").Value = "  +2.01%  "
$ws.Range("D6").Value = "'148.16"
$ws.Range("E6").Value = "  +9.47%  "
$ws.Range("D7").Value = "'3.453.83"
$ws.Range("E7").Value = "  +3.63%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("D9").Value = "'0.473"
$ws.Range("E9").Value = "  +1.28%  "
$ws.Range("D10").Value = "'7.69"
$ws.Range("E10").Value = "  +3.91%  "
$ws.Range("D11").Value = "'0.125"
$ws.Range("E11").Value = "  +2.42%  "
$ws.Range("D12").Value = "'0.387"
$ws.Range("E12").Value = "  +0.96%  "
$ws.Range("D13").Value = "'4.033.69"
$ws.Range("E13").Value = "  +3.00%  "
$ws.Range("D14").Value = "'28.01"
$ws.Range("E14").Value = "  +9.22%  "
$ws.Range("E15").Value = "  -0.26%  "
$ws.Range("D16").Value = "'0.0000176"
$ws.Range("E16").Value = "  +2.51%  "
$ws.Range("D17").Value = "'3.449.65"
$ws.Range("E17").Value = "  +2.96%  "
$ws.Range("D18").Value = "'61.821.27"
$ws.Range("E18").Value = "  +1.16%  "
$ws.Range("D19").Value = "'6.29"
$ws.Range("E19").Value = "  +8.66%  "
$ws.Range("D20").Value = "'14.17"
$ws.Range("E20").Value = "  +2.13%  "
$ws.Range("D21").Value = "'9.41"
$ws.Range("E21").Value = "  +2.39%  "
$ws.Range("D22").Value = "'383.15"
$ws.Range("E22").Value = "  +2.11%  "
$ws.Range("E23").Value = "  +3.17%  "
$ws.Range("D24").Value = "'3.585.13"
$ws.Range("E24").Value = "  +2.53%  "
$ws.Range("D27").Value = "'72.25"
$ws.Range("E27").Value = "  +1.97%  "
$ws.Range("E28").Value = "  +1.01%  "
$ws.Range("D29").Value = "'0.179"
$ws.Range("E29").Value = "  +9.71%  "
$ws.Range("D30").Value = "'7.84"
$ws.Range("E30").Value = "  +5.99%  "
$ws.Range("E31").Value = "  -10.67%  "
$ws.Range("D32").Value = "'0.999"
$ws.Range("E32").Value = "  -0.11%  "
$ws.Range("D33").Value = "'8.20"
$ws.Range("E33").Value = "  +1.44%  "
$ws.Range("E34").Value = "  +2.31%  "
$ws.Range("E35").Value = "  -0.05%  "
$ws.Range("D36").Value = "'23.99"
$ws.Range("E36").Value = "  +2.79%  "
$ws.Range("D37").Value = "'7.03"
$ws.Range("E37").Value = "  +4.51%  "
$ws.Range("D38").Value = "'5.21"
$ws.Range("E38").Value = "  +0.80%  "
$ws.Range("E39").Value = "  +3.10%  "
$ws.Range("D40").Value = "'166.02"
$ws.Range("E40").Value = "  +0.80%  "
$ws.Range("D41").Value = "'0.0784"
$ws.Range("E41").Value = "  +4.14%  "
$ws.Range("B42").Value = "EnergySwap"
$ws.Range("C42").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D42").Value = "'26.09"
$ws.Range("E42").Value = "  +10.07%  "
$ws.Range("B43").Value = "Mantle"
$ws.Range("C43").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D43").Value = "'0.795"
$ws.Range("E43").Value = "  +4.15%  "
$ws.Range("B44").Value = "Stacks"
$ws.Range("C44").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D44").Value = "'1.74"
$ws.Range("E44").Value = "  +3.19%  "
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").Value = "'0.999"
$ws.Range("E45").Value = "  -0.29%  "
$ws.Range("D46").Value = "'42.25"
$ws.Range("E46").Value = "  +2.39%  "
$ws.Range("E47").Value = "  +3.09%  "
$ws.Range("E48").Value = "  -1.37%  "
$ws.Range("D49").Value = "'2.595.90"
$ws.Range("E49").Value = "  +11.00%  "
$ws.Range("D50").Value = "'23.53"
$ws.Range("E50").Value = "  +3.88%  "
$ws.Range("D51").Value = "'6.86"
$ws.Range("E51").Value = "  +1.55%  "
